# Test006 executado e Test007 preparado
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the query table data by the "Weighted" column (E) ascending instead of
# "Friedman Rank_Weighted" (C).
$table = $ws.ListObjects.Item(1)

$table.Sort.SortFields.Clear()
$table.Sort.SortFields.Add($ws.Range("E1:E14"), 0, 1) | Out-Null
$table.Sort.Header = 1
$table.Sort.Apply()

# The sort engine leaves a blank-but-formatted cell behind where a value moved
# out of a column that had holes (e.g. C3, which now has no Friedman
# Rank_Weighted value). Fully clear it so no stray empty cell remains.
$ws.Range("C3").Clear()

# Update the active selection on the sheet to reflect the new cursor position.
$ws.Range("C22").Select()
